$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 115 (shifts the existing row 115 and everything
# below it down by one, growing the used range from A1:R182 to A1:R183).
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row with the new weekly price-report entry.
$ws.Range("A115").Value = 8
$ws.Range("B115").Value = "Terminal La Palmera de La Serena"
$ws.Range("C115").Value = "Coquimbo"
$ws.Range("D115").Value = 45176
$ws.Range("E115").Value = 4
$ws.Range("F115").Value = 100114007
$ws.Range("G115").Value = "Jengibre"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 400
$ws.Range("K115").Value = 16500
$ws.Range("L115").Value = 17000
$ws.Range("M115").Value = 16750
$ws.Range("N115").Value = "$/caja 13 kilos"
$ws.Range("O115").Value = "Perú"
$ws.Range("P115").Value = 1288
$ws.Range("Q115").Value = 13
$ws.Range("R115").Value = "Hortaliza"
